$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 190
$ws.Range("B190").Value = 6810171
$ws.Range("E190").Value = 'KV Kortrijk'
$ws.Range("F190").Value = 'Charleroi'
$ws.Range("G190").Value = 1
$ws.Range("H190").Value = 0
$ws.Range("I190").Value = 'H'
$ws.Range("J190").Value = 3.2
$ws.Range("K190").Value = 3.5
$ws.Range("L190").Value = 2.1
$ws.Range("M190").Value = 3.4
$ws.Range("N190").Value = 3.4
$ws.Range("O190").Value = 2.05
$ws.Range("P190").Value = 0.25
$ws.Range("Q190").Value = 2
$ws.Range("R190").Value = 1.85
$ws.Range("T190").Value = 1.925
$ws.Range("U190").Value = 1.925
$ws.Range("V190").Value = 2.4
$ws.Range("X190").Value = -1
$ws.Range("Y190").Value = 1
$ws.Range("Z190").Value = -1
$ws.Range("AA190").Value = -1
$ws.Range("AB190").Value = 0.925

# Row 191
$ws.Range("B191").Value = 6810174
$ws.Range("E191").Value = 'Westerlo'
$ws.Range("F191").Value = 'OH Leuven'
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 3
$ws.Range("I191").Value = 'A'
$ws.Range("J191").Value = 1.909
$ws.Range("K191").Value = 3.75
$ws.Range("L191").Value = 3.5
$ws.Range("M191").Value = 1.909
$ws.Range("N191").Value = 3.5
$ws.Range("O191").Value = 3.8
$ws.Range("P191").Value = -0.5
$ws.Range("Q191").Value = 1.925
$ws.Range("R191").Value = 1.925
$ws.Range("T191").Value = 1.85
$ws.Range("U191").Value = 2
$ws.Range("V191").Value = -1
$ws.Range("X191").Value = 2.8
$ws.Range("Y191").Value = -1
$ws.Range("Z191").Value = 0.925
$ws.Range("AA191").Value = 0.8500000000000001
$ws.Range("AB191").Value = -1

# Row 241
$ws.Range("B241").Value = 6810219
$ws.Range("E241").Value = 'OH Leuven'
$ws.Range("F241").Value = 'KV Mechelen'
$ws.Range("H241").Value = 0
$ws.Range("I241").Value = 'H'
$ws.Range("J241").Value = 2.8
$ws.Range("K241").Value = 3.5
$ws.Range("L241").Value = 2.375
$ws.Range("M241").Value = 2.7
$ws.Range("N241").Value = 3.5
$ws.Range("O241").Value = 2.45
$ws.Range("P241").Value = 0
$ws.Range("Q241").Value = 2.025
$ws.Range("R241").Value = 1.825
$ws.Range("S241").Value = 2.75
$ws.Range("T241").Value = 1.925
$ws.Range("U241").Value = 1.925
$ws.Range("V241").Value = 1.7
$ws.Range("W241").Value = -1
$ws.Range("Y241").Value = 1.025
$ws.Range("AB241").Value = 0.925

# Row 242
$ws.Range("B242").Value = 6957874
$ws.Range("E242").Value = 'Westerlo'
$ws.Range("F242").Value = 'Genk'
$ws.Range("H242").Value = 1
$ws.Range("I242").Value = 'D'
$ws.Range("J242").Value = 3.6
$ws.Range("K242").Value = 3.8
$ws.Range("L242").Value = 1.909
$ws.Range("M242").Value = 3.6
$ws.Range("N242").Value = 3.75
$ws.Range("O242").Value = 1.95
$ws.Range("P242").Value = 0.5
$ws.Range("Q242").Value = 1.85
$ws.Range("R242").Value = 2
$ws.Range("S242").Value = 3
$ws.Range("T242").Value = 2
$ws.Range("U242").Value = 1.85
$ws.Range("V242").Value = -1
$ws.Range("W242").Value = 2.75
$ws.Range("Y242").Value = 0.8500000000000001
$ws.Range("AB242").Value = 0.8500000000000001

# Row 243
$ws.Range("B243").Value = 6942395
$ws.Range("E243").Value = 'Gent'
$ws.Range("F243").Value = 'Charleroi'
$ws.Range("G243").Value = 5
$ws.Range("J243").Value = 1.571
$ws.Range("K243").Value = 4
$ws.Range("L243").Value = 5.75
$ws.Range("M243").Value = 1.4
$ws.Range("N243").Value = 4.333
$ws.Range("O243").Value = 8
$ws.Range("P243").Value = -1.25
$ws.Range("Q243").Value = 2
$ws.Range("R243").Value = 1.85
$ws.Range("S243").Value = 2.75
$ws.Range("V243").Value = 0.3999999999999999
$ws.Range("Y243").Value = 1

# Row 244
$ws.Range("B244").Value = 6870199
$ws.Range("E244").Value = 'Cercle Brugge'
$ws.Range("F244").Value = 'RWD Molenbeek'
$ws.Range("G244").Value = 4
$ws.Range("J244").Value = 1.363
$ws.Range("K244").Value = 5.5
$ws.Range("L244").Value = 7.5
$ws.Range("M244").Value = 1.3
$ws.Range("N244").Value = 6
$ws.Range("O244").Value = 8.5
$ws.Range("P244").Value = -1.75
$ws.Range("Q244").Value = 2.025
$ws.Range("R244").Value = 1.825
$ws.Range("S244").Value = 3.25
$ws.Range("V244").Value = 0.3
$ws.Range("Y244").Value = 1.025

# Row 276
$ws.Range("B276").Value = 7979470
$ws.Range("E276").Value = 'Westerlo'
$ws.Range("F276").Value = 'OH Leuven'
$ws.Range("G276").Value = 1
$ws.Range("H276").Value = 1
$ws.Range("I276").Value = 'D'
$ws.Range("J276").Value = 2.5
$ws.Range("K276").Value = 3.6
$ws.Range("L276").Value = 2.6
$ws.Range("M276").Value = 2.45
$ws.Range("N276").Value = 3.75
$ws.Range("O276").Value = 2.55
$ws.Range("P276").Value = 0
$ws.Range("Q276").Value = 1.875
$ws.Range("R276").Value = 1.975
$ws.Range("T276").Value = 1.85
$ws.Range("U276").Value = 2
$ws.Range("W276").Value = 2.75
$ws.Range("X276").Value = -1
$ws.Range("Y276").Value = 0
$ws.Range("Z276").Value = 0
$ws.Range("AB276").Value = 1

# Row 277
$ws.Range("B277").Value = 7979346
$ws.Range("E277").Value = 'SintTruidense'
$ws.Range("F277").Value = 'Gent'
$ws.Range("G277").Value = 0
$ws.Range("H277").Value = 2
$ws.Range("I277").Value = 'A'
$ws.Range("J277").Value = 3.6
$ws.Range("K277").Value = 3.5
$ws.Range("L277").Value = 2
$ws.Range("M277").Value = 3.3
$ws.Range("N277").Value = 3.6
$ws.Range("O277").Value = 2.05
$ws.Range("P277").Value = 0.25
$ws.Range("Q277").Value = 2.025
$ws.Range("R277").Value = 1.825
$ws.Range("T277").Value = 1.975
$ws.Range("U277").Value = 1.875
$ws.Range("W277").Value = -1
$ws.Range("X277").Value = 1.05
$ws.Range("Y277").Value = -1
$ws.Range("Z277").Value = 0.825
$ws.Range("AB277").Value = 0.875

# Row 310
$ws.Range("M310").Value = 1.45
$ws.Range("N310").Value = 5
$ws.Range("O310").Value = 5.5
$ws.Range("P310").Value = -1.25
$ws.Range("Q310").Value = 1.95
$ws.Range("R310").Value = 1.9
$ws.Range("T310").Value = 1.85
$ws.Range("U310").Value = 2

# Row 311
$ws.Range("M311").Value = 2
$ws.Range("N311").Value = 3.8
$ws.Range("O311").Value = 3.3
$ws.Range("P311").Value = -0.5
$ws.Range("Q311").Value = 2.025
$ws.Range("R311").Value = 1.825
$ws.Range("T311").Value = 1.875
$ws.Range("U311").Value = 1.975

# Row 312
$ws.Range("N312").Value = 4
$ws.Range("O312").Value = 3.5
$ws.Range("Q312").Value = 1.9
$ws.Range("R312").Value = 1.95
$ws.Range("S312").Value = 3.5
$ws.Range("T312").Value = 2.025
$ws.Range("U312").Value = 1.825

# Row 313
$ws.Range("M313").Value = 3.8
$ws.Range("N313").Value = 3.5
$ws.Range("O313").Value = 1.75
$ws.Range("Q313").Value = 1.8
$ws.Range("R313").Value = 2.05
$ws.Range("T313").Value = 1.95
$ws.Range("U313").Value = 1.9

# Row 314
$ws.Range("M314").Value = 1.55
$ws.Range("N314").Value = 4
$ws.Range("O314").Value = 5
$ws.Range("P314").Value = -1
$ws.Range("Q314").Value = 2
$ws.Range("R314").Value = 1.85
$ws.Range("T314").Value = 1.9
$ws.Range("U314").Value = 1.95

# Row 315
$ws.Range("M315").Value = 1.48
$ws.Range("N315").Value = 4.1
$ws.Range("O315").Value = 5.5
$ws.Range("Q315").Value = 1.85
$ws.Range("R315").Value = 2
$ws.Range("S315").Value = 3.25
$ws.Range("T315").Value = 2
$ws.Range("U315").Value = 1.85
